$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking price values that must stay TEXT (as in the
# original inlineStr cells). Flip the cell to the Text number format before
# assigning so Excel does not silently coerce the string into a Number, then
# clear the now-unneeded explicit format so the cell keeps using the default
# style (matches the untouched sibling cells).
function Set-TextValue($ref, $value) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

Set-TextValue "D2" '243.48'
Set-TextValue "D3" '23.19'
Set-TextValue "D4" '5.410'
Set-TextValue "D5" '0.05968'
Set-TextValue "D6" '3.439'
Set-TextValue "D7" '6.530'
Set-TextValue "D8" '0.8090'
Set-TextValue "D9" '0.9308'
Set-TextValue "D10" '0.1424'
Set-TextValue "D11" '0.07432'
Set-TextValue "D12" '0.03276'
Set-TextValue "D13" '0.03078'
Set-TextValue "D14" '0.09361'
Set-TextValue "D15" '3.860'
Set-TextValue "D16" '0.001570'
Set-TextValue "D17" '0.04712'
Set-TextValue "D18" '0.0005934'
$ws.Range("E18").Value = '17OneONE'
Set-TextValue "D19" '0.005910'
Set-TextValue "D20" '0.001278'
Set-TextValue "D21" '0.004896'
Set-TextValue "D22" '0.00006802'
Set-TextValue "D23" '3.567'
Set-TextValue "D25" '0.3233'
Set-TextValue "D40" '0.03966'
$ws.Range("B41").Value = 'BKEXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
Set-TextValue "D41" '0.1078'
$ws.Range("E41").Value = '40BKEXTokenBKK'
$ws.Range("B42").Value = 'CEJI'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
Set-TextValue "D42" '0.002611'
$ws.Range("E42").Value = '41CEJICEJI'
$ws.Range("B43").Value = 'KickToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
Set-TextValue "D43" '0.003076'
$ws.Range("E43").Value = '42KickTokenKICKWorstin24h'
Set-TextValue "D44" '0.009221'
$ws.Range("E44").Value = '43LocalTradersLCTBestin24h'
Set-TextValue "D45" '0.00005138'
Set-TextValue "D46" '0.00000000751'
Set-TextValue "D47" '0.7005'
Set-TextValue "D48" '0.002411'
Set-TextValue "D49" '0.00002102'
